# Update "想去人数" (interested-count) figures in column F across all four sheets
# to reflect the latest generated output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 360
$ws.Range("F5").Value = 1120
$ws.Range("F8").Value = 888
$ws.Range("F9").Value = 1607
$ws.Range("F10").Value = 6054
$ws.Range("F11").Value = 110
$ws.Range("F12").Value = 1737
$ws.Range("F13").Value = 440
$ws.Range("F14").Value = 5888
$ws.Range("F15").Value = 113
$ws.Range("F19").Value = 1640
$ws.Range("F23").Value = 1361
$ws.Range("F24").Value = 720
$ws.Range("F25").Value = 235
$ws.Range("F30").Value = 3855

# --- Sheet "演出" ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 307
$ws.Range("F5").Value = 158

# --- Sheet "本地生活" ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 9493
$ws.Range("F3").Value = 2226
$ws.Range("F4").Value = 614
$ws.Range("F5").Value = 175

# --- Sheet "全部类型" ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 9493
$ws.Range("F3").Value = 2226
$ws.Range("F4").Value = 614
$ws.Range("F5").Value = 360
$ws.Range("F7").Value = 1120
$ws.Range("F11").Value = 307
$ws.Range("F12").Value = 888
$ws.Range("F13").Value = 175
$ws.Range("F14").Value = 1607
$ws.Range("F15").Value = 6054
$ws.Range("F16").Value = 110
$ws.Range("F17").Value = 1737
$ws.Range("F20").Value = 440
$ws.Range("F23").Value = 5888
$ws.Range("F24").Value = 113
$ws.Range("F28").Value = 1640
$ws.Range("F32").Value = 1361
$ws.Range("F33").Value = 720
$ws.Range("F35").Value = 235
$ws.Range("F45").Value = 3855
